$d = $word.ActiveDocument

# 1. Promote the title paragraph from Heading1 to Heading2
$d.Paragraphs.Item(1).Range.Style = "Heading2"

# 2. The report table now only needs Metric / Value - drop the
#    "Data Lineage" column and widen the remaining two to fill the space.
$t = $d.Tables.Item(1)
$t.Columns.Item(3).Delete()
$t.Columns.Item(1).Width = 216
$t.Columns.Item(2).Width = 216

# 3. Replace the terse snake_case metric names with the friendlier,
#    human-readable labels (the Value column is left untouched).
$t.Cell(2,1).Range.Text = "'adjudicator employment term, hearing centre, do not use reason' table missing columns"
$t.Cell(3,1).Range.Text = "'adjudicator employment term, hearing centre, do not use reason' table data type mismatch count"
$t.Cell(4,1).Range.Text = "'judicial officer history, users' table missing columns"
$t.Cell(5,1).Range.Text = "'judicial officer history, users' table data type mismatch count"
$t.Cell(6,1).Range.Text = "'other centre hearing centre' table missing columns"
$t.Cell(7,1).Range.Text = "'other centre hearing centre' table data type mismatch count"
$t.Cell(8,1).Range.Text = "'adjudicator role' table missing columns"
$t.Cell(9,1).Range.Text = "'adjudicator role' table data type mismatch count"
$t.Cell(10,1).Range.Text = "'adjudicator employment term, hearing centre, do not use reason' table exists"
$t.Cell(11,1).Range.Text = "'adjudicator employment term, hearing centre, do not use reason' 'adjudicator id' table null count"
$t.Cell(12,1).Range.Text = "'judicial officer history, users' table exists"
$t.Cell(13,1).Range.Text = "'judicial officer history, users' 'adjudicator id' table null count"
$t.Cell(14,1).Range.Text = "'other centre hearing centre' table exists"
$t.Cell(15,1).Range.Text = "'other centre hearing centre' 'adjudicator id' table null count"
$t.Cell(16,1).Range.Text = "'adjudicator role' table exists"
$t.Cell(17,1).Range.Text = "'adjudicator role' 'adjudicator id' table null count"
